$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '67.938.48'
$ws.Range('E2').Value = '  -1.50%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.835.47'
$ws.Range('E3').Value = '  -2.39%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '597.51'
$ws.Range('E5').Value = '  -1.20%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '166.62'
$ws.Range('E6').Value = '  +0.84%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '3.828.84'
$ws.Range('E7').Value = '  -2.44%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.527'
$ws.Range('E9').Value = '  -0.95%  '
$ws.Range('E10').Value = '  -1.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.30'
$ws.Range('E11').Value = '  -1.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.460'
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000248'
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '37.13'
$ws.Range('E14').Value = '  -0.40%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '4.492.38'
$ws.Range('E15').Value = '  -1.99%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.833.59'
$ws.Range('E16').Value = '  -2.91%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '68.240.25'
$ws.Range('E17').Value = '  -1.22%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.53'
$ws.Range('E18').Value = '  +0.66%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '18.10'
$ws.Range('E19').Value = '  +5.45%  '
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.72'
$ws.Range('E21').Value = '  -3.85%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '468.68'
$ws.Range('E22').Value = '  -4.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.731'
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000159'
$ws.Range('E24').Value = '  -4.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '84.31'
$ws.Range('E25').Value = '  -0.28%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.22'
$ws.Range('E26').Value = '  -1.91%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.29'
$ws.Range('E27').Value = '  +1.23%  '
$ws.Range('B28').Value = 'Dai'
$ws.Range('C28').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.00'
$ws.Range('E28').Value = '  -0.08%  '
$ws.Range('B29').Value = 'RenderToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.99'
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.93'
$ws.Range('E30').Value = '  -0.69%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.996.98'
$ws.Range('E31').Value = '  -2.03%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '7.72'
$ws.Range('E32').Value = '  -1.99%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.30'
$ws.Range('E33').Value = '  -4.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '31.00'
$ws.Range('E34').Value = '  -4.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.820.75'
$ws.Range('E35').Value = '  -1.40%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.105'
$ws.Range('E36').Value = '  -1.96%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.139'
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.92'
$ws.Range('E38').Value = '  -0.31%  '
$ws.Range('B39').Value = 'Mantle'
$ws.Range('C39').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.00'
$ws.Range('E39').Value = '  -3.84%  '
$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.27'
$ws.Range('E40').Value = '  +8.33%  '
$ws.Range('E41').Value = '  +0.01%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.313'
$ws.Range('E42').Value = '  -2.79%  '
$ws.Range('B43').Value = 'Stacks'
$ws.Range('C43').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.00'
$ws.Range('E43').Value = '  -0.37%  '
$ws.Range('B44').Value = 'Bittensor'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '426.81'
$ws.Range('E44').Value = '  -3.39%  '
$ws.Range('E46').Value = '  -2.08%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '8.58'
$ws.Range('E47').Value = '  +1.29%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000271'
$ws.Range('E48').Value = '  +13.05%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '142.37'
$ws.Range('E49').Value = '  +0.19%  '
$ws.Range('E50').Value = '  +0.23%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '39.06'
$ws.Range('E51').Value = '  -0.47%  '
